# feat: add 2022-Q4 data
#
# Before: 总计, 2022-Q3, 2022-Q2
# After:  总计, 2022-Q4 (new), 2022-Q3, 2022-Q2
#
# The new "2022-Q4" sheet re-uses the same fund roster (codes/names) as the
# sheet that used to be called "2022-Q3", just with refreshed numbers, so we
# duplicate that sheet and patch the changed cells. The "总计" (totals) sheet
# gets a new row for the 2022-Q4 quarter, inserted above the existing rows.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)   # currently named "2022-Q3"

# Helper: write a string value into a cell while forcing a genuine text
# (non-numeric) cell type, matching the source data where numeric-looking
# figures ("10.01", "0.1311", ...) are stored as text, not numbers.
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# ---------------------------------------------------------------------
# 1. Duplicate the "2022-Q3" sheet to become the new "2022-Q4" sheet,
#    placed right after "总计" (i.e. right before the old "2022-Q3").
# ---------------------------------------------------------------------
$q3Sheet.Copy($q3Sheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# Patch the new Q4 sheet's numbers (fund codes/names stay the same).
Set-TextValue $q4Sheet.Range("D2") "10.01"
Set-TextValue $q4Sheet.Range("E2") "93.81"
Set-TextValue $q4Sheet.Range("F2") "1.31"
Set-TextValue $q4Sheet.Range("G2") "0.1311"
$q4Sheet.Range("H2").Value = 1

Set-TextValue $q4Sheet.Range("D3") "0.94"
Set-TextValue $q4Sheet.Range("E3") "97.66"
Set-TextValue $q4Sheet.Range("F3") "1.37"
Set-TextValue $q4Sheet.Range("G3") "0.0129"
$q4Sheet.Range("H3").Value = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: shift the existing two rows down
#    and insert the new 2022-Q4 row on top of them.
# ---------------------------------------------------------------------

# Row 4 (new) <- old row 3 data (2022-Q2 / 0.15), copying A3's style for A4.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.15

# Row 3 <- old row 2 data (2022-Q3 / 0.13).
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.13

# Row 2 <- brand-new 2022-Q4 data.
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.14

# Restore "总计" as the active sheet, matching the workbook's original
# bookViews (activeTab="0"), which the diff leaves untouched.
$totalSheet.Activate() | Out-Null
